# BRVM Recommandations – mise à jour automatique (GitHub Actions)
# Refreshes the "Recommandations" sector/stock table (rows 2-44) and the
# "Top_YTD" progression table (rows 2-11) with the latest computed values,
# re-ranks rows by "Variation Totale (%)" and grows the used range from
# A1:G42 to A1:G44 (two new rows: AIR LIQUIDE CI (SIVC) and ORAGROUP TOGO (ORGT)).

$wb     = $excel.ActiveWorkbook
$wsReco = $wb.Worksheets.Item("Recommandations")
$wsYtd  = $wb.Worksheets.Item("Top_YTD")

# --- "Recommandations" sheet: rows 2-44, columns A-G ---
$wsReco.Range("A2").Value = "BRVM - SERVICES PUBLICS"
$wsReco.Range("B2").Value = 0
$wsReco.Range("C2").Value = 8
$wsReco.Range("D2").Value = 3354.62
$wsReco.Range("E2").Value = 105.89
$wsReco.Range("F2").Value = "🟡 Observer"
$wsReco.Range("G2").Value = "➖ Neutre"

$wsReco.Range("A3").Value = "NEI-CEDA CI"
$wsReco.Range("B3").Value = 0
$wsReco.Range("C3").Value = 4
$wsReco.Range("D3").Value = 2735
$wsReco.Range("E3").Value = 680
$wsReco.Range("F3").Value = "🟡 Observer"
$wsReco.Range("G3").Value = "➖ Neutre"

$wsReco.Range("A4").Value = "AIR LIQUIDE CI"
$wsReco.Range("B4").Value = 0
$wsReco.Range("C4").Value = 4
$wsReco.Range("D4").Value = 2645
$wsReco.Range("E4").Value = 650
$wsReco.Range("F4").Value = "🟡 Observer"
$wsReco.Range("G4").Value = "➖ Neutre"

$wsReco.Range("A5").Value = "BRVM - AUTRES SECTEURS"
$wsReco.Range("B5").Value = 0
$wsReco.Range("C5").Value = 4
$wsReco.Range("D5").Value = 2332.95
$wsReco.Range("E5").Value = 609.68
$wsReco.Range("F5").Value = "🟡 Observer"
$wsReco.Range("G5").Value = "➖ Neutre"

$wsReco.Range("A6").Value = "BRVM - DISTRIBUTION"
$wsReco.Range("B6").Value = 0
$wsReco.Range("C6").Value = 4
$wsReco.Range("D6").Value = 1812.9
$wsReco.Range("E6").Value = 477.82
$wsReco.Range("F6").Value = "🟡 Observer"
$wsReco.Range("G6").Value = "➖ Neutre"

$wsReco.Range("A7").Value = "BRVM - AGRICULTURE"
$wsReco.Range("B7").Value = 0
$wsReco.Range("C7").Value = 4
$wsReco.Range("D7").Value = 1475.69
$wsReco.Range("E7").Value = 386.43
$wsReco.Range("F7").Value = "🟡 Observer"
$wsReco.Range("G7").Value = "➖ Neutre"

$wsReco.Range("A8").Value = "BRVM - TRANSPORT"
$wsReco.Range("B8").Value = 0
$wsReco.Range("C8").Value = 4
$wsReco.Range("D8").Value = 1439.27
$wsReco.Range("E8").Value = 364.71
$wsReco.Range("F8").Value = "🟡 Observer"
$wsReco.Range("G8").Value = "➖ Neutre"

$wsReco.Range("A9").Value = "BRVM - CONSOMMATION DISCRETIONNAIRE"
$wsReco.Range("B9").Value = 0
$wsReco.Range("C9").Value = 4
$wsReco.Range("D9").Value = 593.84
$wsReco.Range("E9").Value = 161.03
$wsReco.Range("F9").Value = "🟡 Observer"
$wsReco.Range("G9").Value = "➖ Neutre"

$wsReco.Range("A10").Value = "BRVM-PRESTIGE"
$wsReco.Range("B10").Value = 0
$wsReco.Range("C10").Value = 4
$wsReco.Range("D10").Value = 548.16
$wsReco.Range("E10").Value = 136.79
$wsReco.Range("F10").Value = "🟡 Observer"
$wsReco.Range("G10").Value = "➖ Neutre"

$wsReco.Range("A11").Value = "BRVM - FINANCES"
$wsReco.Range("B11").Value = 0
$wsReco.Range("C11").Value = 4
$wsReco.Range("D11").Value = 542.34
$wsReco.Range("E11").Value = 135.15
$wsReco.Range("F11").Value = "🟡 Observer"
$wsReco.Range("G11").Value = "➖ Neutre"

$wsReco.Range("A12").Value = "BRVM - SERVICES FINANCIERS"
$wsReco.Range("B12").Value = 0
$wsReco.Range("C12").Value = 4
$wsReco.Range("D12").Value = 533.01
$wsReco.Range("E12").Value = 132.82
$wsReco.Range("F12").Value = "🟡 Observer"
$wsReco.Range("G12").Value = "➖ Neutre"

$wsReco.Range("A13").Value = "BRVM - INDUSTRIELS"
$wsReco.Range("B13").Value = 0
$wsReco.Range("C13").Value = 4
$wsReco.Range("D13").Value = 505.01
$wsReco.Range("E13").Value = 131.68
$wsReco.Range("F13").Value = "🟡 Observer"
$wsReco.Range("G13").Value = "➖ Neutre"

$wsReco.Range("A14").Value = "BRVM - ENERGIE"
$wsReco.Range("B14").Value = 0
$wsReco.Range("C14").Value = 4
$wsReco.Range("D14").Value = 435.9
$wsReco.Range("E14").Value = 110.39
$wsReco.Range("F14").Value = "🟡 Observer"
$wsReco.Range("G14").Value = "➖ Neutre"

$wsReco.Range("A15").Value = "BRVM - TELECOMMUNICATIONS"
$wsReco.Range("B15").Value = 0
$wsReco.Range("C15").Value = 4
$wsReco.Range("D15").Value = 382.6
$wsReco.Range("E15").Value = 96.62
$wsReco.Range("F15").Value = "🟡 Observer"
$wsReco.Range("G15").Value = "➖ Neutre"

$wsReco.Range("A16").Value = "BRVM - INDUSTRIE                        (**)"
$wsReco.Range("B16").Value = 0
$wsReco.Range("C16").Value = 1
$wsReco.Range("D16").Value = 212.57
$wsReco.Range("E16").Value = 212.57
$wsReco.Range("F16").Value = "🟡 Observer"
$wsReco.Range("G16").Value = "➖ Neutre"

$wsReco.Range("A17").Value = "BRVM - INDUSTRIE                  (**)"
$wsReco.Range("B17").Value = 0
$wsReco.Range("C17").Value = 1
$wsReco.Range("D17").Value = 207.68
$wsReco.Range("E17").Value = 207.68
$wsReco.Range("F17").Value = "🟡 Observer"
$wsReco.Range("G17").Value = "➖ Neutre"

$wsReco.Range("A18").Value = "BRVM-PRINCIPAL                           (**)"
$wsReco.Range("B18").Value = 0
$wsReco.Range("C18").Value = 1
$wsReco.Range("D18").Value = 199.2
$wsReco.Range("E18").Value = 199.2
$wsReco.Range("F18").Value = "🟡 Observer"
$wsReco.Range("G18").Value = "➖ Neutre"

$wsReco.Range("A19").Value = "BRVM-PRINCIPAL                    (**)"
$wsReco.Range("B19").Value = 0
$wsReco.Range("C19").Value = 1
$wsReco.Range("D19").Value = 193.83
$wsReco.Range("E19").Value = 193.83
$wsReco.Range("F19").Value = "🟡 Observer"
$wsReco.Range("G19").Value = "➖ Neutre"

$wsReco.Range("A20").Value = "BRVM - CONSOMMATION DE BASE             (**)"
$wsReco.Range("B20").Value = 0
$wsReco.Range("C20").Value = 1
$wsReco.Range("D20").Value = 192.99
$wsReco.Range("E20").Value = 192.99
$wsReco.Range("F20").Value = "🟡 Observer"
$wsReco.Range("G20").Value = "➖ Neutre"

$wsReco.Range("A21").Value = "BRVM - CONSOMMATION DE BASE         (**)"
$wsReco.Range("B21").Value = 0
$wsReco.Range("C21").Value = 1
$wsReco.Range("D21").Value = 187.07
$wsReco.Range("E21").Value = 187.07
$wsReco.Range("F21").Value = "🟡 Observer"
$wsReco.Range("G21").Value = "➖ Neutre"

$wsReco.Range("A22").Value = "TRACTAFRIC MOTORS CI (PRSC)"
$wsReco.Range("B22").Value = 4
$wsReco.Range("C22").Value = 0
$wsReco.Range("D22").Value = 29.59
$wsReco.Range("E22").Value = 7.4
$wsReco.Range("F22").Value = "🟢 Achat"
$wsReco.Range("G22").Value = "✅ Renforcer"

$wsReco.Range("A23").Value = "FILTISAC CI (FTSC)"
$wsReco.Range("B23").Value = 4
$wsReco.Range("C23").Value = 0
$wsReco.Range("D23").Value = 29.5
$wsReco.Range("E23").Value = 7.37
$wsReco.Range("F23").Value = "🟢 Achat"
$wsReco.Range("G23").Value = "✅ Renforcer"

$wsReco.Range("A24").Value = "CFAO MOTORS CI (CFAC)"
$wsReco.Range("B24").Value = 3
$wsReco.Range("C24").Value = 0
$wsReco.Range("D24").Value = 17.98
$wsReco.Range("E24").Value = 7.35
$wsReco.Range("F24").Value = "🟢 Achat"
$wsReco.Range("G24").Value = "✅ Renforcer"

$wsReco.Range("A25").Value = "BERNABE CI (BNBC)"
$wsReco.Range("B25").Value = 2
$wsReco.Range("C25").Value = 0
$wsReco.Range("D25").Value = 14.01
$wsReco.Range("E25").Value = 7.5
$wsReco.Range("F25").Value = "🟡 Observer"
$wsReco.Range("G25").Value = "➖ Neutre"

$wsReco.Range("A26").Value = "SAPH CI (SPHC)"
$wsReco.Range("B26").Value = 2
$wsReco.Range("C26").Value = 0
$wsReco.Range("D26").Value = 13.39
$wsReco.Range("E26").Value = 7.45
$wsReco.Range("F26").Value = "🟡 Observer"
$wsReco.Range("G26").Value = "➖ Neutre"

$wsReco.Range("A27").Value = "SICABLE CI (CABC)"
$wsReco.Range("B27").Value = 1
$wsReco.Range("C27").Value = 0
$wsReco.Range("D27").Value = 7.33
$wsReco.Range("E27").Value = 7.33
$wsReco.Range("F27").Value = "🟡 Observer"
$wsReco.Range("G27").Value = "➖ Neutre"

$wsReco.Range("A28").Value = "SUCRIVOIRE (SCRC)"
$wsReco.Range("B28").Value = 1
$wsReco.Range("C28").Value = 1
$wsReco.Range("D28").Value = 3.35
$wsReco.Range("E28").Value = 6.78
$wsReco.Range("F28").Value = "🟡 Observer"
$wsReco.Range("G28").Value = "➖ Neutre"

$wsReco.Range("A29").Value = "VIVO ENERGY CI (SHEC)"
$wsReco.Range("B29").Value = 1
$wsReco.Range("C29").Value = 0
$wsReco.Range("D29").Value = 2.92
$wsReco.Range("E29").Value = 2.92
$wsReco.Range("F29").Value = "🟡 Observer"
$wsReco.Range("G29").Value = "➖ Neutre"

$wsReco.Range("A30").Value = "SETAO CI (STAC)"
$wsReco.Range("B30").Value = 1
$wsReco.Range("C30").Value = 2
$wsReco.Range("D30").Value = 2.51
$wsReco.Range("E30").Value = 7.32
$wsReco.Range("F30").Value = "🟡 Observer"
$wsReco.Range("G30").Value = "➖ Neutre"

$wsReco.Range("A31").Value = "TOTAL"
$wsReco.Range("B31").Value = 0
$wsReco.Range("C31").Value = 5
$wsReco.Range("D31").Value = 0
$wsReco.Range("E31").Value = 0
$wsReco.Range("F31").Value = "🟡 Observer"
$wsReco.Range("G31").Value = "➖ Neutre"

$wsReco.Range("A32").Value = "UNIWAX CI (UNXC)"
$wsReco.Range("B32").Value = 1
$wsReco.Range("C32").Value = 1
$wsReco.Range("D32").Value = -0.25
$wsReco.Range("E32").Value = 7.23
$wsReco.Range("F32").Value = "🟡 Observer"
$wsReco.Range("G32").Value = "👀 À surveiller"

$wsReco.Range("A33").Value = "BANK OF AFRICA ML (BOAM)"
$wsReco.Range("B33").Value = 0
$wsReco.Range("C33").Value = 1
$wsReco.Range("D33").Value = -1.15
$wsReco.Range("E33").Value = -1.15
$wsReco.Range("F33").Value = "🟡 Observer"
$wsReco.Range("G33").Value = "➖ Neutre"

$wsReco.Range("A34").Value = "SOGB CI (SOGC)"
$wsReco.Range("B34").Value = 0
$wsReco.Range("C34").Value = 1
$wsReco.Range("D34").Value = -1.18
$wsReco.Range("E34").Value = -1.18
$wsReco.Range("F34").Value = "🟡 Observer"
$wsReco.Range("G34").Value = "➖ Neutre"

$wsReco.Range("A35").Value = "AFRICA GLOBAL LOGISTICS CI (SDSC)"
$wsReco.Range("B35").Value = 0
$wsReco.Range("C35").Value = 1
$wsReco.Range("D35").Value = -1.68
$wsReco.Range("E35").Value = -1.68
$wsReco.Range("F35").Value = "🟡 Observer"
$wsReco.Range("G35").Value = "➖ Neutre"

$wsReco.Range("A36").Value = "SOCIETE GENERALE COTE D'IVOIRE (SGBC)"
$wsReco.Range("B36").Value = 0
$wsReco.Range("C36").Value = 1
$wsReco.Range("D36").Value = -1.85
$wsReco.Range("E36").Value = -1.85
$wsReco.Range("F36").Value = "🟡 Observer"
$wsReco.Range("G36").Value = "➖ Neutre"

$wsReco.Range("A37").Value = "ONATEL BF (ONTBF)"
$wsReco.Range("B37").Value = 0
$wsReco.Range("C37").Value = 1
$wsReco.Range("D37").Value = -1.9
$wsReco.Range("E37").Value = -1.9
$wsReco.Range("F37").Value = "🟡 Observer"
$wsReco.Range("G37").Value = "➖ Neutre"

$wsReco.Range("A38").Value = "ECOBANK COTE D''IVOIRE (ECOC)"
$wsReco.Range("B38").Value = 0
$wsReco.Range("C38").Value = 1
$wsReco.Range("D38").Value = -2.23
$wsReco.Range("E38").Value = -2.23
$wsReco.Range("F38").Value = "🟡 Observer"
$wsReco.Range("G38").Value = "➖ Neutre"

$wsReco.Range("A39").Value = "SOLIBRA CI (SLBC)"
$wsReco.Range("B39").Value = 0
$wsReco.Range("C39").Value = 1
$wsReco.Range("D39").Value = -2.75
$wsReco.Range("E39").Value = -2.75
$wsReco.Range("F39").Value = "🟡 Observer"
$wsReco.Range("G39").Value = "➖ Neutre"

$wsReco.Range("A40").Value = "ORANGE COTE D'IVOIRE (ORAC)"
$wsReco.Range("B40").Value = 0
$wsReco.Range("C40").Value = 1
$wsReco.Range("D40").Value = -3.34
$wsReco.Range("E40").Value = -3.34
$wsReco.Range("F40").Value = "🟡 Observer"
$wsReco.Range("G40").Value = "➖ Neutre"

$wsReco.Range("A41").Value = "NSIA BANQUE COTE D'IVOIRE (NSBC)"
$wsReco.Range("B41").Value = 0
$wsReco.Range("C41").Value = 2
$wsReco.Range("D41").Value = -3.84
$wsReco.Range("E41").Value = -1.83
$wsReco.Range("F41").Value = "🟡 Observer"
$wsReco.Range("G41").Value = "➖ Neutre"

$wsReco.Range("A42").Value = "NEI-CEDA CI (NEIC)"
$wsReco.Range("B42").Value = 0
$wsReco.Range("C42").Value = 2
$wsReco.Range("D42").Value = -6.59
$wsReco.Range("E42").Value = -3.65
$wsReco.Range("F42").Value = "🟡 Observer"
$wsReco.Range("G42").Value = "➖ Neutre"

$wsReco.Range("A43").Value = "AIR LIQUIDE CI (SIVC)"
$wsReco.Range("B43").Value = 0
$wsReco.Range("C43").Value = 1
$wsReco.Range("D43").Value = -7.35
$wsReco.Range("E43").Value = -7.35
$wsReco.Range("F43").Value = "🟡 Observer"
$wsReco.Range("G43").Value = "➖ Neutre"

$wsReco.Range("A44").Value = "ORAGROUP TOGO (ORGT)"
$wsReco.Range("B44").Value = 0
$wsReco.Range("C44").Value = 1
$wsReco.Range("D44").Value = -7.38
$wsReco.Range("E44").Value = -7.38
$wsReco.Range("F44").Value = "🟡 Observer"
$wsReco.Range("G44").Value = "➖ Neutre"

# --- "Top_YTD" sheet: refreshed progression figures (row 8 unchanged) ---
$wsYtd.Range("B2").Value = 8588223.41
$wsYtd.Range("B3").Value = 377199.26
$wsYtd.Range("B4").Value = 334958.75
$wsYtd.Range("B5").Value = 217487.55
$wsYtd.Range("B6").Value = 93377.78
$wsYtd.Range("B7").Value = 48205.73
$wsYtd.Range("B9").Value = 3700.15
$wsYtd.Range("B10").Value = 3057.08
$wsYtd.Range("B11").Value = 2980.27

